# edit.ps1 — apply the commit's changes via PowerPoint COM-interop
#
# 1. Slide 5's table switches from the deck's custom "Table_0" style
#    ({93D193A7-909B-4580-B894-6E06F777E87F}) to the built-in table
#    style {B99D3B4A-BC67-4C75-9B92-C53A121C19F5}.
# 2. The presentation's theme palette is swapped from the "Integral"
#    ("Red Violet") scheme to the stock "Office Theme" scheme.

$p = $ppt.ActivePresentation

# --- helper: pack R,G,B (0-255 each) into the COM/OLE RGB() int (0x00BBGGRR) ---
function ToComRgb([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# --- 1. Re-style the table on slide 5 --------------------------------------
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{B99D3B4A-BC67-4C75-9B92-C53A121C19F5}")
    }
}

# --- 2. Swap the theme color scheme to the stock Office palette ------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (in that order)
$officeColors = @(
    @(0x00, 0x00, 0x00),  # dk1
    @(0xFF, 0xFF, 0xFF),  # lt1
    @(0x44, 0x54, 0x6A),  # dk2
    @(0xE7, 0xE6, 0xE6),  # lt2
    @(0x5B, 0x9B, 0xD5),  # accent1
    @(0xED, 0x7D, 0x31),  # accent2
    @(0xA5, 0xA5, 0xA5),  # accent3
    @(0xFF, 0xC0, 0x00),  # accent4
    @(0x44, 0x72, 0xC4),  # accent5
    @(0x70, 0xAD, 0x47),  # accent6
    @(0x05, 0x63, 0xC1),  # hlink
    @(0x95, 0x4F, 0x72)   # folHlink
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $c = $officeColors[$i - 1]
    $themeColors.Item($i).RGB = ToComRgb $c[0] $c[1] $c[2]
}
